$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'FAPs'
$ws.Cells.Item(2,2).Value = 'Fgf21'
$ws.Cells.Item(2,3).Value = 'Fgfr4'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.2030856666666666
$ws.Cells.Item(2,8).Value = 0.6092569999999999
$ws.Cells.Item(2,9).Value = 0.3531826360773729
$ws.Cells.Item(2,10).Value = 0.3531826360773729
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.903185
$ws.Cells.Item(2,14).Value = 2.709555
$ws.Cells.Item(2,15).Value = 0.03154869388788047
$ws.Cells.Item(2,16).Value = 0.03154869388788046
$ws.Cells.Item(2,17).Value = 0.1834239278483333
$ws.Cells.Item(2,18).Value = 1.650815350635
$ws.Cells.Item(2,19).Value = 0.01114245087211973
$ws.Cells.Item(2,20).Value = 0.01114245087211972

# Row 3
$ws.Cells.Item(3,1).Value = 'FAPs'
$ws.Cells.Item(3,2).Value = 'Fgf21'
$ws.Cells.Item(3,3).Value = 'Fgfr4'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.2030856666666666
$ws.Cells.Item(3,8).Value = 0.6092569999999999
$ws.Cells.Item(3,9).Value = 0.3531826360773729
$ws.Cells.Item(3,10).Value = 0.3531826360773729
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.110372
$ws.Cells.Item(3,14).Value = 0.331116
$ws.Cells.Item(3,15).Value = 0.003855347953955327
$ws.Cells.Item(3,16).Value = 0.003855347953955326
$ws.Cells.Item(3,17).Value = 0.02241497120133333
$ws.Cells.Item(3,18).Value = 0.201734740812
$ws.Cells.Item(3,19).Value = 0.001361641953373448
$ws.Cells.Item(3,20).Value = 0.001361641953373448

# Row 4
$ws.Cells.Item(4,1).Value = 'FAPs'
$ws.Cells.Item(4,2).Value = 'Fgf21'
$ws.Cells.Item(4,3).Value = 'Fgfr4'
$ws.Cells.Item(4,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.2030856666666666
$ws.Cells.Item(4,8).Value = 0.6092569999999999
$ws.Cells.Item(4,9).Value = 0.3531826360773729
$ws.Cells.Item(4,10).Value = 0.3531826360773729
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.03076233333333333
$ws.Cells.Item(4,14).Value = 0.09228699999999999
$ws.Cells.Item(4,15).Value = 0.001074543352259254
$ws.Cells.Item(4,16).Value = 0.001074543352259254
$ws.Cells.Item(4,17).Value = 0.006247388973222221
$ws.Cells.Item(4,18).Value = 0.05622650075899999
$ws.Cells.Item(4,19).Value = 0.0003795100537303405
$ws.Cells.Item(4,20).Value = 0.0003795100537303405

# Row 5
$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'Fgf21'
$ws.Cells.Item(5,3).Value = 'Fgfr4'
$ws.Cells.Item(5,4).Value = 'MuSCs'
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.2030856666666666
$ws.Cells.Item(5,8).Value = 0.6092569999999999
$ws.Cells.Item(5,9).Value = 0.3531826360773729
$ws.Cells.Item(5,10).Value = 0.3531826360773729
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 27.52907633333334
$ws.Cells.Item(5,14).Value = 82.58722900000001
$ws.Cells.Item(5,15).Value = 0.9616041035407232
$ws.Cells.Item(5,16).Value = 0.9616041035407231
$ws.Cells.Item(5,17).Value = 5.590760819872555
$ws.Cells.Item(5,18).Value = 50.316847378853
$ws.Cells.Item(5,19).Value = 0.3396218721513317
$ws.Cells.Item(5,20).Value = 0.3396218721513316

# Row 6
$ws.Cells.Item(6,1).Value = 'FAPs'
$ws.Cells.Item(6,2).Value = 'Fgf21'
$ws.Cells.Item(6,3).Value = 'Fgfr4'
$ws.Cells.Item(6,4).Value = 'Resolving-Mac'
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.2030856666666666
$ws.Cells.Item(6,8).Value = 0.6092569999999999
$ws.Cells.Item(6,9).Value = 0.3531826360773729
$ws.Cells.Item(6,10).Value = 0.3531826360773729
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.05488933333333334
$ws.Cells.Item(6,14).Value = 0.164668
$ws.Cells.Item(6,15).Value = 0.001917311265181737
$ws.Cells.Item(6,16).Value = 0.001917311265181736
$ws.Cells.Item(6,17).Value = 0.01114723685288889
$ws.Cells.Item(6,18).Value = 0.100325131676
$ws.Cells.Item(6,19).Value = 0.0006771610468177286
$ws.Cells.Item(6,20).Value = 0.0006771610468177286

# Row 7
$ws.Cells.Item(7,1).Value = 'MuSCs'
$ws.Cells.Item(7,2).Value = 'Fgf21'
$ws.Cells.Item(7,3).Value = 'Fgfr4'
$ws.Cells.Item(7,4).Value = 'ECs'
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.3719303333333333
$ws.Cells.Item(7,8).Value = 1.115791
$ws.Cells.Item(7,9).Value = 0.6468173639226271
$ws.Cells.Item(7,10).Value = 0.6468173639226271
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.903185
$ws.Cells.Item(7,14).Value = 2.709555
$ws.Cells.Item(7,15).Value = 0.03154869388788047
$ws.Cells.Item(7,16).Value = 0.03154869388788046
$ws.Cells.Item(7,17).Value = 0.3359218981116667
$ws.Cells.Item(7,18).Value = 3.023297083005
$ws.Cells.Item(7,19).Value = 0.02040624301576074
$ws.Cells.Item(7,20).Value = 0.02040624301576074

# Row 8
$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Fgf21'
$ws.Cells.Item(8,3).Value = 'Fgfr4'
$ws.Cells.Item(8,4).Value = 'FAPs'
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.3719303333333333
$ws.Cells.Item(8,8).Value = 1.115791
$ws.Cells.Item(8,9).Value = 0.6468173639226271
$ws.Cells.Item(8,10).Value = 0.6468173639226271
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.110372
$ws.Cells.Item(8,14).Value = 0.331116
$ws.Cells.Item(8,15).Value = 0.003855347953955327
$ws.Cells.Item(8,16).Value = 0.003855347953955326
$ws.Cells.Item(8,17).Value = 0.04105069475066666
$ws.Cells.Item(8,18).Value = 0.369456252756
$ws.Cells.Item(8,19).Value = 0.002493706000581878
$ws.Cells.Item(8,20).Value = 0.002493706000581878

# Row 9
$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Fgf21'
$ws.Cells.Item(9,3).Value = 'Fgfr4'
$ws.Cells.Item(9,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.3719303333333333
$ws.Cells.Item(9,8).Value = 1.115791
$ws.Cells.Item(9,9).Value = 0.6468173639226271
$ws.Cells.Item(9,10).Value = 0.6468173639226271
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.03076233333333333
$ws.Cells.Item(9,14).Value = 0.09228699999999999
$ws.Cells.Item(9,15).Value = 0.001074543352259254
$ws.Cells.Item(9,16).Value = 0.001074543352259254
$ws.Cells.Item(9,17).Value = 0.01144144489077778
$ws.Cells.Item(9,18).Value = 0.102973004017
$ws.Cells.Item(9,19).Value = 0.0006950332985289139
$ws.Cells.Item(9,20).Value = 0.0006950332985289139

# Row 10
$ws.Cells.Item(10,1).Value = 'MuSCs'
$ws.Cells.Item(10,2).Value = 'Fgf21'
$ws.Cells.Item(10,3).Value = 'Fgfr4'
$ws.Cells.Item(10,4).Value = 'MuSCs'
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.3719303333333333
$ws.Cells.Item(10,8).Value = 1.115791
$ws.Cells.Item(10,9).Value = 0.6468173639226271
$ws.Cells.Item(10,10).Value = 0.6468173639226271
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 27.52907633333334
$ws.Cells.Item(10,14).Value = 82.58722900000001
$ws.Cells.Item(10,15).Value = 0.9616041035407232
$ws.Cells.Item(10,16).Value = 0.9616041035407231
$ws.Cells.Item(10,17).Value = 10.23889853701544
$ws.Cells.Item(10,18).Value = 92.15008683313901
$ws.Cells.Item(10,19).Value = 0.6219822313893916
$ws.Cells.Item(10,20).Value = 0.6219822313893916

# Row 11
$ws.Cells.Item(11,1).Value = 'MuSCs'
$ws.Cells.Item(11,2).Value = 'Fgf21'
$ws.Cells.Item(11,3).Value = 'Fgfr4'
$ws.Cells.Item(11,4).Value = 'Resolving-Mac'
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.3719303333333333
$ws.Cells.Item(11,8).Value = 1.115791
$ws.Cells.Item(11,9).Value = 0.6468173639226271
$ws.Cells.Item(11,10).Value = 0.6468173639226271
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.05488933333333334
$ws.Cells.Item(11,14).Value = 0.164668
$ws.Cells.Item(11,15).Value = 0.001917311265181737
$ws.Cells.Item(11,16).Value = 0.001917311265181736
$ws.Cells.Item(11,17).Value = 0.02041500804311111
$ws.Cells.Item(11,18).Value = 0.183735072388
$ws.Cells.Item(11,19).Value = 0.001240150218364008
$ws.Cells.Item(11,20).Value = 0.001240150218364008
